$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.089.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.960.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "380.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.546"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.94%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.588"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.32%  "

$ws.Range("E11").Value = "  -0.67%  "

$ws.Range("E12").Value = "  +1.95%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.424.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +69.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.951.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.152.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.80%  "

$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.44%  "

$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0965"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +16.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.34%  "

$ws.Range("E26").Value = "  -2.25%  "

$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("E28").Value = "  -0.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.110"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.06%  "

$ws.Range("E33").Value = "  +5.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "34.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0438"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.20%  "

$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.25%  "

$ws.Range("E39").Value = "  +2.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.22%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.53%  "

$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.62%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.278"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.065.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.41%  "

$ws.Range("E48").Value = "  -1.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.90%  "

$ws.Range("E50").Value = "  -8.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.00%  "
